# Fruta / hortaliza, semanal
# Reorders rows 2-20 by re-shuffling the D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg) values
# according to the new weekly ordering, while leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: target row -> source row (where the new values for D,M,N,O,P,S come from)
$mapping = @{
    2  = 8
    3  = 18
    4  = 5
    5  = 3
    6  = 17
    7  = 11
    8  = 7
    9  = 6
    10 = 16
    11 = 14
    12 = 20
    13 = 10
    14 = 19
    15 = 2
    16 = 9
    17 = 12
    18 = 13
    19 = 4
    20 = 15
}

# Capture the original values for the columns that move, before any writes happen.
# NOTE: use Value2 for reads - Value's getter in this runtime does not return the
# underlying scalar, only Value2 (and Text) do. The Value *setter* works fine.
$orig = @{}
for ($r = 2; $r -le 20; $r++) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $vals = $orig[$src]
    $ws.Cells.Item($r, 4).Value = $vals.D
    $ws.Cells.Item($r, 13).Value = $vals.M
    $ws.Cells.Item($r, 14).Value = $vals.N
    $ws.Cells.Item($r, 15).Value = $vals.O
    $ws.Cells.Item($r, 16).Value = $vals.P
    $ws.Cells.Item($r, 19).Value = $vals.S
}
